$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 851.5833
$ws.Range("I11").Value = 851.5833
$ws.Range("K11").Value = 851.5833
$ws.Range("M11").Value = -711.5833
$ws.Range("H28").Value = 3273.9092
$ws.Range("I28").Value = 2719.3215
$ws.Range("J28").Value = 6379.6
$ws.Range("K28").Value = 2719.3215
$ws.Range("L28").Value = 6379.6
$ws.Range("M28").Value = -2234.3215
$ws.Range("N28").Value = -7349.6
$ws.Range("H70").Value = 17089.428
$ws.Range("J70").Value = 19104.334
$ws.Range("L70").Value = 57313.00199999999
$ws.Range("N70").Value = -57853.00199999999
$ws.Range("H73").Value = 17089.428
$ws.Range("J73").Value = 19104.334
$ws.Range("L73").Value = 57313.00199999999
$ws.Range("N73").Value = -59185.00199999999
$ws.Range("H86").Value = 3565.7856
$ws.Range("I86").Value = 3174.6365
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3174.6365
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2051.6365
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 3565.7856
$ws.Range("I89").Value = 3174.6365
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 15873.1825
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -10257.1825
$ws.Range("N89").Value = -36232
$ws.Range("H92").Value = 260.57895
$ws.Range("I92").Value = 230.92308
$ws.Range("K92").Value = 230.92308
$ws.Range("M92").Value = 1017.07692
$ws.Range("H106").Value = 5688082.5
$ws.Range("I106").Value = 7582449.5
$ws.Range("J106").Value = 4982.5
$ws.Range("K106").Value = 7582449.5
$ws.Range("L106").Value = 4982.5
$ws.Range("M106").Value = -7581818.5
$ws.Range("N106").Value = -6244.5
$ws.Range("H107").Value = 476.6154
$ws.Range("I107").Value = 502
$ws.Range("K107").Value = 502
$ws.Range("M107").Value = 1418
$ws.Range("H111").Value = 1900
$ws.Range("I111").Value = 1900
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 5700
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -2633
$ws.Range("N111").ClearContents()
$ws.Range("H132").Value = 4980.4062
$ws.Range("I132").Value = 5319.5557
$ws.Range("J132").Value = 3149
$ws.Range("K132").Value = 15958.6671
$ws.Range("L132").Value = 9447
$ws.Range("M132").Value = -13428.6671
$ws.Range("N132").Value = -14507
$ws.Range("H134").Value = 31928.938
$ws.Range("J134").Value = 31928.938
$ws.Range("L134").Value = 31928.938
$ws.Range("N134").Value = -42068.93799999999
$ws.Range("H138").Value = 2774.978
$ws.Range("I138").Value = 4412.3184
$ws.Range("J138").Value = 2252.9275
$ws.Range("K138").Value = 13236.9552
$ws.Range("L138").Value = 6758.782499999999
$ws.Range("M138").Value = -8096.9552
$ws.Range("N138").Value = -17038.7825

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 696.65515
$ws.Range("I2").Value = 598.2308
$ws.Range("J2").Value = 1549.6666
$ws.Range("K2").Value = 598.2308
$ws.Range("L2").Value = 1549.6666
$ws.Range("M2").Value = -485.2308
$ws.Range("N2").Value = -1775.6666
$ws.Range("H32").Value = 4447.8096
$ws.Range("I32").Value = 4188.9033
$ws.Range("K32").Value = 4188.9033
$ws.Range("M32").Value = -3901.9033
$ws.Range("H38").Value = 2506.3333
$ws.Range("I38").Value = 2506.3333
$ws.Range("K38").Value = 2506.3333
$ws.Range("M38").Value = -2039.3333
$ws.Range("H39").Value = 25000
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H45").Value = 3669.3333
$ws.Range("I45").Value = 2789.3333
$ws.Range("J45").Value = 4021.3333
$ws.Range("K45").Value = 2789.3333
$ws.Range("L45").Value = 4021.3333
$ws.Range("M45").Value = -2412.3333
$ws.Range("N45").Value = -4775.3333
$ws.Range("H61").Value = 10627.2
$ws.Range("I61").Value = 2874
$ws.Range("K61").Value = 2874
$ws.Range("M61").Value = -2662
$ws.Range("H63").Value = 790
$ws.Range("I63").Value = 790
$ws.Range("K63").Value = 790
$ws.Range("M63").Value = -104
$ws.Range("H66").Value = 790
$ws.Range("I66").Value = 790
$ws.Range("K66").Value = 3950
$ws.Range("M66").Value = -518
$ws.Range("H74").Value = 1474.3
$ws.Range("I74").Value = 957.6667
$ws.Range("J74").Value = 2249.25
$ws.Range("K74").Value = 957.6667
$ws.Range("L74").Value = 2249.25
$ws.Range("M74").Value = -83.66669999999999
$ws.Range("N74").Value = -3997.25
$ws.Range("H77").Value = 1474.3
$ws.Range("I77").Value = 957.6667
$ws.Range("J77").Value = 2249.25
$ws.Range("K77").Value = 4788.3335
$ws.Range("L77").Value = 11246.25
$ws.Range("M77").Value = -420.3334999999997
$ws.Range("N77").Value = -19982.25
$ws.Range("H88").Value = 6034.385
$ws.Range("I88").Value = 1609.375
$ws.Range("J88").Value = 8001.0557
$ws.Range("K88").Value = 1609.375
$ws.Range("L88").Value = 8001.0557
$ws.Range("M88").Value = -1203.375
$ws.Range("N88").Value = -8813.0557
$ws.Range("H91").Value = 6034.385
$ws.Range("I91").Value = 1609.375
$ws.Range("J91").Value = 8001.0557
$ws.Range("K91").Value = 1609.375
$ws.Range("L91").Value = 8001.0557
$ws.Range("M91").Value = -205.375
$ws.Range("N91").Value = -10809.0557
$ws.Range("H97").Value = 2295.3914
$ws.Range("I97").Value = 496.42856
$ws.Range("K97").Value = 496.42856
$ws.Range("M97").Value = -0.4285600000000045
$ws.Range("H110").Value = 339
$ws.Range("I110").Value = 337.5
$ws.Range("K110").Value = 337.5
$ws.Range("M110").Value = 1707.5
$ws.Range("H116").Value = 696.65515
$ws.Range("I116").Value = 598.2308
$ws.Range("J116").Value = 1549.6666
$ws.Range("K116").Value = 598.2308
$ws.Range("L116").Value = 1549.6666
$ws.Range("M116").Value = 1695.7692
$ws.Range("N116").Value = -6137.6666
$ws.Range("H122").Value = 1503.4762
$ws.Range("I122").Value = 1159.25
$ws.Range("J122").Value = 2605
$ws.Range("K122").Value = 3477.75
$ws.Range("L122").Value = 7815
$ws.Range("M122").Value = -1027.75
$ws.Range("N122").Value = -12715
$ws.Range("H136").Value = 10627.2
$ws.Range("I136").Value = 2874
$ws.Range("K136").Value = 8622
$ws.Range("M136").Value = -6072

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 696.65515
$ws.Range("I3").Value = 598.2308
$ws.Range("J3").Value = 1549.6666
$ws.Range("K3").Value = 598.2308
$ws.Range("L3").Value = 1549.6666
$ws.Range("M3").Value = -484.2308
$ws.Range("N3").Value = -1777.6666
$ws.Range("H6").Value = 80460.08
$ws.Range("J6").Value = 81993.664
$ws.Range("L6").Value = 81993.664
$ws.Range("N6").Value = -82219.664
$ws.Range("H107").Value = 999.2857
$ws.Range("I107").Value = 999.2857
$ws.Range("K107").Value = 999.2857
$ws.Range("M107").Value = 920.7143
$ws.Range("H132").Value = 95999
$ws.Range("J132").Value = 95999
$ws.Range("L132").Value = 95999
$ws.Range("N132").Value = -106119
$ws.Range("H133").Value = 60326.332
$ws.Range("J133").Value = 60326.332
$ws.Range("L133").Value = 60326.332
$ws.Range("N133").Value = -70446.332
$ws.Range("H134").Value = 9284.187
$ws.Range("I134").Value = 8748.634
$ws.Range("J134").Value = 10520.077
$ws.Range("K134").Value = 26245.902
$ws.Range("L134").Value = 31560.231
$ws.Range("M134").Value = -23710.902
$ws.Range("N134").Value = -36630.231
$ws.Range("H141").Value = 66249.5
$ws.Range("I141").Value = 60000
$ws.Range("J141").Value = 68332.664
$ws.Range("K141").Value = 60000
$ws.Range("L141").Value = 68332.664
$ws.Range("M141").Value = -54820
$ws.Range("N141").Value = -78692.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1767.3158
$ws.Range("I16").Value = 1161.875
$ws.Range("J16").Value = 4996.3335
$ws.Range("K16").Value = 1161.875
$ws.Range("L16").Value = 4996.3335
$ws.Range("M16").Value = -874.875
$ws.Range("N16").Value = -5570.3335
$ws.Range("H31").Value = 2083.919
$ws.Range("I31").Value = 1726.5714
$ws.Range("J31").Value = 2552.9375
$ws.Range("K31").Value = 1726.5714
$ws.Range("L31").Value = 2552.9375
$ws.Range("M31").Value = -1431.5714
$ws.Range("N31").Value = -3142.9375
$ws.Range("H34").Value = 2083.919
$ws.Range("I34").Value = 1726.5714
$ws.Range("J34").Value = 2552.9375
$ws.Range("K34").Value = 1726.5714
$ws.Range("L34").Value = 2552.9375
$ws.Range("M34").Value = -1524.5714
$ws.Range("N34").Value = -2956.9375
$ws.Range("H58").Value = 2563.5454
$ws.Range("I58").Value = 2265.5625
$ws.Range("J58").Value = 3358.1667
$ws.Range("K58").Value = 2265.5625
$ws.Range("L58").Value = 3358.1667
$ws.Range("M58").Value = -2062.5625
$ws.Range("N58").Value = -3764.1667
$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41498
$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -127488
$ws.Range("H74").Value = 39998
$ws.Range("J74").Value = 39998
$ws.Range("L74").Value = 39998
$ws.Range("N74").Value = -41746
$ws.Range("H77").Value = 39998
$ws.Range("J77").Value = 39998
$ws.Range("L77").Value = 119994
$ws.Range("N77").Value = -128730
$ws.Range("H94").Value = 2433.5715
$ws.Range("I94").Value = 2236.1667
$ws.Range("K94").Value = 2236.1667
$ws.Range("M94").Value = -1785.1667
$ws.Range("H99").Value = 2322.7778
$ws.Range("I99").Value = 2329.3572
$ws.Range("K99").Value = 2329.3572
$ws.Range("M99").Value = -831.3571999999999
$ws.Range("H113").Value = 1767.3158
$ws.Range("I113").Value = 1161.875
$ws.Range("J113").Value = 4996.3335
$ws.Range("K113").Value = 1161.875
$ws.Range("L113").Value = 4996.3335
$ws.Range("M113").Value = 1008.125
$ws.Range("N113").Value = -9336.3335
$ws.Range("H122").Value = 3441.4092
$ws.Range("J122").Value = 3631.7778
$ws.Range("L122").Value = 10895.3334
$ws.Range("N122").Value = -15795.3334
$ws.Range("H126").Value = 2322.7778
$ws.Range("I126").Value = 2329.3572
$ws.Range("K126").Value = 6988.071599999999
$ws.Range("M126").Value = -4518.071599999999
$ws.Range("H132").Value = 3368.3333
$ws.Range("I132").Value = 3011.3125
$ws.Range("J132").Value = 4510.8
$ws.Range("K132").Value = 9033.9375
$ws.Range("L132").Value = 13532.4
$ws.Range("M132").Value = -6503.9375
$ws.Range("N132").Value = -18592.4
$ws.Range("H134").Value = 1971.4706
$ws.Range("I134").Value = 1359.5555
$ws.Range("K134").Value = 4078.6665
$ws.Range("M134").Value = -1543.6665
$ws.Range("H136").Value = 2563.5454
$ws.Range("I136").Value = 2265.5625
$ws.Range("J136").Value = 3358.1667
$ws.Range("K136").Value = 6796.6875
$ws.Range("L136").Value = 10074.5001
$ws.Range("M136").Value = -4246.6875
$ws.Range("N136").Value = -15174.5001
$ws.Range("H141").Value = 349878.56
$ws.Range("J141").Value = 457630
$ws.Range("L141").Value = 457630
$ws.Range("N141").Value = -467990

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 556.4167
$ws.Range("I5").Value = 527.7
$ws.Range("K5").Value = 1583.1
$ws.Range("M5").Value = -1471.1
$ws.Range("H8").Value = 793.2857
$ws.Range("I8").Value = 793.2857
$ws.Range("K8").Value = 2379.8571
$ws.Range("M8").Value = -2240.8571
$ws.Range("H34").Value = 2504.6667
$ws.Range("J34").Value = 10861.6
$ws.Range("L34").Value = 32584.8
$ws.Range("N34").Value = -32752.8
$ws.Range("H37").Value = 126613.336
$ws.Range("J37").Value = 126613.336
$ws.Range("L37").Value = 379840.008
$ws.Range("N37").Value = -380064.008
$ws.Range("H39").Value = 10371.143
$ws.Range("I39").Value = 2300
$ws.Range("J39").Value = 13599.6
$ws.Range("K39").Value = 6900
$ws.Range("L39").Value = 40798.8
$ws.Range("M39").Value = -6606
$ws.Range("N39").Value = -41386.8
$ws.Range("H50").Value = 1514.7
$ws.Range("I50").Value = 2869.75
$ws.Range("J50").Value = 611.3333
$ws.Range("K50").Value = 8609.25
$ws.Range("L50").Value = 1833.9999
$ws.Range("M50").Value = -8128.25
$ws.Range("N50").Value = -2795.9999
$ws.Range("H53").Value = 1514.7
$ws.Range("I53").Value = 2869.75
$ws.Range("J53").Value = 611.3333
$ws.Range("K53").Value = 8609.25
$ws.Range("L53").Value = 1833.9999
$ws.Range("M53").Value = -8128.25
$ws.Range("N53").Value = -2795.9999
$ws.Range("H56").Value = 7024.595
$ws.Range("I56").Value = 7024.595
$ws.Range("K56").Value = 7024.595
$ws.Range("M56").Value = -6494.595
$ws.Range("H125").Value = 17500
$ws.Range("I125").Value = 8750
$ws.Range("K125").Value = 26250
$ws.Range("M125").Value = -21330
$ws.Range("H131").Value = 19324514
$ws.Range("I131").Value = 6536649.5
$ws.Range("J131").Value = 55556796
$ws.Range("K131").Value = 19609948.5
$ws.Range("L131").Value = 166670388
$ws.Range("M131").Value = -19604908.5
$ws.Range("N131").Value = -166680468
$ws.Range("H132").Value = 4379.636
$ws.Range("I132").Value = 1848.5
$ws.Range("J132").Value = 4942.1113
$ws.Range("K132").Value = 16636.5
$ws.Range("L132").Value = 44479.00169999999
$ws.Range("M132").Value = -14106.5
$ws.Range("N132").Value = -49539.00169999999
$ws.Range("H135").Value = 556.4167
$ws.Range("I135").Value = 527.7
$ws.Range("K135").Value = 4749.3
$ws.Range("M135").Value = -2214.3
$ws.Range("H136").Value = 7281
$ws.Range("I136").Value = 5074.364
$ws.Range("K136").Value = 15223.092
$ws.Range("M136").Value = -10123.092
$ws.Range("H140").Value = 5962649
$ws.Range("I140").Value = 14707782
$ws.Range("J140").Value = 15959.2
$ws.Range("K140").Value = 44123346
$ws.Range("L140").Value = 47877.60000000001
$ws.Range("M140").Value = -44118166
$ws.Range("N140").Value = -58237.60000000001
$ws.Range("H141").Value = 32105.709
$ws.Range("I141").Value = 8030.778
$ws.Range("J141").Value = 41954.547
$ws.Range("K141").Value = 24092.334
$ws.Range("L141").Value = 125863.641
$ws.Range("M141").Value = -18912.334
$ws.Range("N141").Value = -136223.641

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10516.643
$ws.Range("I70").Value = 11248.6
$ws.Range("K70").Value = 11248.6
$ws.Range("M70").Value = -10978.6
$ws.Range("H73").Value = 10516.643
$ws.Range("I73").Value = 11248.6
$ws.Range("K73").Value = 11248.6
$ws.Range("M73").Value = -10312.6
$ws.Range("H80").Value = 4026.125
$ws.Range("I80").Value = 3652.5557
$ws.Range("J80").Value = 4506.4287
$ws.Range("K80").Value = 3652.5557
$ws.Range("L80").Value = 4506.4287
$ws.Range("M80").Value = -2654.5557
$ws.Range("N80").Value = -6502.4287
$ws.Range("H83").Value = 4026.125
$ws.Range("I83").Value = 3652.5557
$ws.Range("J83").Value = 4506.4287
$ws.Range("K83").Value = 18262.7785
$ws.Range("L83").Value = 22532.1435
$ws.Range("M83").Value = -13270.7785
$ws.Range("N83").Value = -32516.1435
$ws.Range("H97").Value = 805.56757
$ws.Range("I97").Value = 582.56
$ws.Range("J97").Value = 1270.1666
$ws.Range("K97").Value = 582.56
$ws.Range("L97").Value = 1270.1666
$ws.Range("M97").Value = -86.55999999999995
$ws.Range("N97").Value = -2262.1666
$ws.Range("H113").Value = 1441.3125
$ws.Range("I113").Value = 1069.25
$ws.Range("J113").Value = 1565.3334
$ws.Range("K113").Value = 1069.25
$ws.Range("L113").Value = 1565.3334
$ws.Range("M113").Value = 1100.75
$ws.Range("N113").Value = -5905.3334
$ws.Range("H122").Value = 5980.737
$ws.Range("I122").Value = 5622.909
$ws.Range("J122").Value = 6472.75
$ws.Range("K122").Value = 16868.727
$ws.Range("L122").Value = 19418.25
$ws.Range("M122").Value = -14418.727
$ws.Range("N122").Value = -24318.25
$ws.Range("H132").Value = 3822.25
$ws.Range("I132").Value = 2939.4614
$ws.Range("J132").Value = 5461.7144
$ws.Range("K132").Value = 8818.3842
$ws.Range("L132").Value = 16385.1432
$ws.Range("M132").Value = -6288.3842
$ws.Range("N132").Value = -21445.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 623.05554
$ws.Range("I16").Value = 559.5714
$ws.Range("J16").Value = 845.25
$ws.Range("K16").Value = 559.5714
$ws.Range("L16").Value = 845.25
$ws.Range("M16").Value = -389.5714
$ws.Range("N16").Value = -1185.25
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H61").Value = 1744.6842
$ws.Range("I61").Value = 1702.7778
$ws.Range("J61").Value = 2499
$ws.Range("K61").Value = 1702.7778
$ws.Range("L61").Value = 2499
$ws.Range("M61").Value = -1500.7778
$ws.Range("N61").Value = -2903
$ws.Range("H68").Value = 2029.88
$ws.Range("I68").Value = 1979.619
$ws.Range("K68").Value = 1979.619
$ws.Range("M68").Value = -1230.619
$ws.Range("H70").Value = 28932
$ws.Range("I70").Value = 14997
$ws.Range("J70").Value = 35899.5
$ws.Range("K70").Value = 14997
$ws.Range("L70").Value = 35899.5
$ws.Range("M70").Value = -14727
$ws.Range("N70").Value = -36439.5
$ws.Range("H71").Value = 2029.88
$ws.Range("I71").Value = 1979.619
$ws.Range("K71").Value = 9898.095
$ws.Range("M71").Value = -6154.094999999999
$ws.Range("H73").Value = 28932
$ws.Range("I73").Value = 14997
$ws.Range("J73").Value = 35899.5
$ws.Range("K73").Value = 14997
$ws.Range("L73").Value = 35899.5
$ws.Range("M73").Value = -14061
$ws.Range("N73").Value = -37771.5
$ws.Range("H82").Value = 2775.0908
$ws.Range("I82").Value = 917.44446
$ws.Range("J82").Value = 4061.1538
$ws.Range("K82").Value = 917.44446
$ws.Range("L82").Value = 4061.1538
$ws.Range("M82").Value = -556.44446
$ws.Range("N82").Value = -4783.1538
$ws.Range("H85").Value = 2775.0908
$ws.Range("I85").Value = 917.44446
$ws.Range("J85").Value = 4061.1538
$ws.Range("K85").Value = 917.44446
$ws.Range("L85").Value = 4061.1538
$ws.Range("M85").Value = 330.55554
$ws.Range("N85").Value = -6557.1538
$ws.Range("H100").Value = 4399.3335
$ws.Range("I100").Value = 3749
$ws.Range("J100").Value = 4724.5
$ws.Range("K100").Value = 3749
$ws.Range("L100").Value = 4724.5
$ws.Range("M100").Value = -3208
$ws.Range("N100").Value = -5806.5
$ws.Range("H113").Value = 1744.6842
$ws.Range("I113").Value = 1702.7778
$ws.Range("J113").Value = 2499
$ws.Range("K113").Value = 1702.7778
$ws.Range("L113").Value = 2499
$ws.Range("M113").Value = 467.2221999999999
$ws.Range("N113").Value = -6839
$ws.Range("H122").Value = 3093.75
$ws.Range("J122").Value = 3093.75
$ws.Range("L122").Value = 9281.25
$ws.Range("N122").Value = -14181.25
$ws.Range("H132").Value = 3981.3057
$ws.Range("I132").Value = 3815.5
$ws.Range("J132").Value = 4188.5625
$ws.Range("K132").Value = 11446.5
$ws.Range("L132").Value = 12565.6875
$ws.Range("M132").Value = -8916.5
$ws.Range("N132").Value = -17625.6875
$ws.Range("H136").Value = 2933.0435
$ws.Range("I136").Value = 2655.862
$ws.Range("J136").Value = 3405.8823
$ws.Range("K136").Value = 7967.586
$ws.Range("L136").Value = 10217.6469
$ws.Range("M136").Value = -5417.586
$ws.Range("N136").Value = -15317.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 35000
$ws.Range("J28").Value = 35000
$ws.Range("L28").Value = 35000
$ws.Range("N28").Value = -35696
$ws.Range("H33").Value = 21446.5
$ws.Range("J33").Value = 27595.334
$ws.Range("L33").Value = 27595.334
$ws.Range("N33").Value = -28095.334
$ws.Range("H36").Value = 21446.5
$ws.Range("J36").Value = 27595.334
$ws.Range("L36").Value = 27595.334
$ws.Range("N36").Value = -28095.334
$ws.Range("H44").Value = 44475
$ws.Range("J44").Value = 44633.332
$ws.Range("L44").Value = 44633.332
$ws.Range("N44").Value = -45741.332
$ws.Range("H96").Value = 1214.2858
$ws.Range("I96").Value = 1291.6666
$ws.Range("J96").Value = 750
$ws.Range("K96").Value = 1291.6666
$ws.Range("L96").Value = 750
$ws.Range("M96").Value = 81.33339999999998
$ws.Range("N96").Value = -3496
$ws.Range("H107").Value = 8117.143
$ws.Range("I107").Value = 13033.5
$ws.Range("K107").Value = 39100.5
$ws.Range("M107").Value = -37180.5
$ws.Range("H113").Value = 526.8571
$ws.Range("I113").Value = 774.75
$ws.Range("J113").Value = 196.33333
$ws.Range("K113").Value = 2324.25
$ws.Range("L113").Value = 588.99999
$ws.Range("M113").Value = -154.25
$ws.Range("N113").Value = -4928.99999
$ws.Range("H126").Value = 3312.25
$ws.Range("I126").Value = 2625
$ws.Range("K126").Value = 7875
$ws.Range("M126").Value = -5405
$ws.Range("H132").Value = 2020.3448
$ws.Range("I132").Value = 1800.88
$ws.Range("K132").Value = 5402.64
$ws.Range("M132").Value = -2872.64
